# Auto-generated Excel COM-interop script applying the Aegis_Profits market-data refresh
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (hunk 0)
$ws.Range("H12").Value = 100310
$ws.Range("I12").Value = 416.66666
$ws.Range("J12").Value = 250150
$ws.Range("K12").Value = 416.66666
$ws.Range("L12").Value = 250150
$ws.Range("M12").Value = -246.66666
$ws.Range("N12").Value = -250490

# Row 137 (hunk 1)
$ws.Range("H137").Value = 1965.6154
$ws.Range("I137").Value = 2081.25
$ws.Range("J137").Value = 1780.6
$ws.Range("K137").Value = 6243.75
$ws.Range("L137").Value = 5341.799999999999
$ws.Range("M137").Value = -3693.75
$ws.Range("N137").Value = -10441.8

# Row 138 (hunk 2)
$ws.Range("H138").Value = 3793.3242
$ws.Range("I138").Value = 4427.2856
$ws.Range("J138").Value = 3727.0896
$ws.Range("K138").Value = 13281.8568
$ws.Range("L138").Value = 11181.2688
$ws.Range("M138").Value = -8141.856800000001
$ws.Range("N138").Value = -21461.2688

# Row 141 (hunk 3)
$ws.Range("H141").Value = 2962.3157
$ws.Range("I141").Value = 2565.6
$ws.Range("J141").Value = 4450
$ws.Range("K141").Value = 7696.799999999999
$ws.Range("L141").Value = 13350
$ws.Range("M141").Value = -2516.799999999999
$ws.Range("N141").Value = -23710

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 4)
$ws.Range("H32").Value = 40432.492
$ws.Range("I32").Value = 14434.48
$ws.Range("J32").Value = 226132.58
$ws.Range("K32").Value = 14434.48
$ws.Range("L32").Value = 226132.58
$ws.Range("M32").Value = -14147.48
$ws.Range("N32").Value = -226706.58

# Row 61 (hunk 5)
$ws.Range("H61").Value = 1864.721
$ws.Range("I61").Value = 1646.6316
$ws.Range("J61").Value = 2037.375
$ws.Range("K61").Value = 1646.6316
$ws.Range("L61").Value = 2037.375
$ws.Range("M61").Value = -1434.6316
$ws.Range("N61").Value = -2461.375

# Row 101 (hunk 6)
$ws.Range("H101").Value = 28301
$ws.Range("J101").Value = 28301
$ws.Range("L101").Value = 28301
$ws.Range("N101").Value = -34791

# Row 132 (hunk 7)
$ws.Range("H132").Value = 10534.807
$ws.Range("I132").Value = 12561.94
$ws.Range("J132").Value = 2088.4167
$ws.Range("K132").Value = 37685.82
$ws.Range("L132").Value = 6265.250100000001
$ws.Range("M132").Value = -35155.82
$ws.Range("N132").Value = -11325.2501

# Row 136 (hunk 8)
$ws.Range("H136").Value = 1864.721
$ws.Range("I136").Value = 1646.6316
$ws.Range("J136").Value = 2037.375
$ws.Range("K136").Value = 4939.8948
$ws.Range("L136").Value = 6112.125
$ws.Range("M136").Value = -2389.8948
$ws.Range("N136").Value = -11212.125

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (hunk 9)
$ws.Range("H22").Value = 296.66666
$ws.Range("I22").Value = 296.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 296.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -123.66666
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (hunk 10)
$ws.Range("H4").Value = 4000500
$ws.Range("I4").Value = 5000750
$ws.Range("J4").Value = 2000000
$ws.Range("K4").Value = 5000750
$ws.Range("L4").Value = 2000000
$ws.Range("M4").Value = -5000638
$ws.Range("N4").Value = -2000224

# Row 22 (hunk 11)
$ws.Range("H22").Value = 946
$ws.Range("I22").Value = 228
$ws.Range("K22").Value = 228
$ws.Range("M22").Value = 122

# Row 31 (hunk 12)
$ws.Range("H31").Value = 45081.91
$ws.Range("I31").Value = 1435.4286
$ws.Range("J31").Value = 75634.45
$ws.Range("K31").Value = 1435.4286
$ws.Range("L31").Value = 75634.45
$ws.Range("M31").Value = -1140.4286
$ws.Range("N31").Value = -76224.45

# Row 34 (hunk 13)
$ws.Range("H34").Value = 45081.91
$ws.Range("I34").Value = 1435.4286
$ws.Range("J34").Value = 75634.45
$ws.Range("K34").Value = 1435.4286
$ws.Range("L34").Value = 75634.45
$ws.Range("M34").Value = -1233.4286
$ws.Range("N34").Value = -76038.45

# Row 58 (hunk 14)
$ws.Range("H58").Value = 1529.7675
$ws.Range("I58").Value = 1345.6072
$ws.Range("J58").Value = 1873.5333
$ws.Range("K58").Value = 1345.6072
$ws.Range("L58").Value = 1873.5333
$ws.Range("M58").Value = -1142.6072
$ws.Range("N58").Value = -2279.5333

# Row 122 (hunk 15)
$ws.Range("H122").Value = 1390.3334
$ws.Range("J122").Value = 1390.3334
$ws.Range("L122").Value = 4171.0002
$ws.Range("N122").Value = -9071.0002

# Row 134 (hunk 16)
$ws.Range("H134").Value = 1375.8096
$ws.Range("I134").Value = 777.6923
$ws.Range("J134").Value = 2347.75
$ws.Range("K134").Value = 2333.0769
$ws.Range("L134").Value = 7043.25
$ws.Range("M134").Value = 201.9231
$ws.Range("N134").Value = -12113.25

# Row 136 (hunk 17)
$ws.Range("H136").Value = 1529.7675
$ws.Range("I136").Value = 1345.6072
$ws.Range("J136").Value = 1873.5333
$ws.Range("K136").Value = 4036.8216
$ws.Range("L136").Value = 5620.5999
$ws.Range("M136").Value = -1486.8216
$ws.Range("N136").Value = -10720.5999

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (hunk 18)
$ws.Range("H4").Value = 750
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 3000
$ws.Range("N4").Value = -3224

# Row 59 (hunk 19)
$ws.Range("H59").Value = 916.6667
$ws.Range("I59").Value = 750
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 2250
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = -1710
$ws.Range("N59").Value = -4080

# Row 114 (hunk 20)
$ws.Range("H114").Value = 1166.6923
$ws.Range("I114").Value = 628.375
$ws.Range("J114").Value = 2028
$ws.Range("K114").Value = 1885.125
$ws.Range("L114").Value = 6084
$ws.Range("M114").Value = 1368.875
$ws.Range("N114").Value = -12592

# Row 120 (hunk 21)
$ws.Range("H120").Value = 502015
$ws.Range("I120").Value = 502015
$ws.Range("K120").Value = 1506045
$ws.Range("M120").Value = -1501207

# Row 122 (hunk 22)
$ws.Range("H122").Value = 788.5625
$ws.Range("I122").Value = 474.14285
$ws.Range("J122").Value = 1033.1111
$ws.Range("K122").Value = 4267.28565
$ws.Range("L122").Value = 9297.999900000001
$ws.Range("M122").Value = -1817.28565
$ws.Range("N122").Value = -14197.9999

# Row 141 (hunk 23)
$ws.Range("H141").Value = 3990
$ws.Range("I141").Value = 4231.4287
$ws.Range("J141").Value = 2300
$ws.Range("K141").Value = 12694.2861
$ws.Range("L141").Value = 6900
$ws.Range("M141").Value = -7514.286100000001
$ws.Range("N141").Value = -17260

$ws = $wb.Worksheets.Item("GSM")
# Row 5 (hunk 24)
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 80 (hunk 25)
$ws.Range("H80").Value = 200444900
$ws.Range("J80").Value = 19500
$ws.Range("L80").Value = 19500
$ws.Range("N80").Value = -21496

# Row 83 (hunk 26)
$ws.Range("H83").Value = 200444900
$ws.Range("J83").Value = 19500
$ws.Range("L83").Value = 97500
$ws.Range("N83").Value = -107484

# Row 101 (hunk 27)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 132 (hunk 28)
$ws.Range("H132").Value = 2334.132
$ws.Range("I132").Value = 1774.4736
$ws.Range("J132").Value = 3751.9333
$ws.Range("K132").Value = 5323.4208
$ws.Range("L132").Value = 11255.7999
$ws.Range("M132").Value = -2793.4208
$ws.Range("N132").Value = -16315.7999

$ws = $wb.Worksheets.Item("LTW")
# Row 2 (hunk 29)
$ws.Range("H2").Value = 289285.72
$ws.Range("J2").Value = 9200
$ws.Range("L2").Value = 9200
$ws.Range("N2").Value = -9424

# Row 101 (hunk 30)
$ws.Range("H101").Value = 26403
$ws.Range("J101").Value = 26403
$ws.Range("L101").Value = 26403
$ws.Range("N101").Value = -32893

# Row 103 (hunk 31)
$ws.Range("H103").Value = 46187
$ws.Range("J103").Value = 46187
$ws.Range("L103").Value = 46187
$ws.Range("N103").Value = -48531

# Row 104 (hunk 32)
$ws.Range("H104").Value = 18603.334
$ws.Range("J104").Value = 18603.334
$ws.Range("L104").Value = 18603.334
$ws.Range("N104").Value = -25591.334

# Row 132 (hunk 33)
$ws.Range("H132").Value = 4329.6
$ws.Range("I132").Value = 5610.769
$ws.Range("J132").Value = 2941.6667
$ws.Range("K132").Value = 16832.307
$ws.Range("L132").Value = 8825.000100000001
$ws.Range("M132").Value = -14302.307
$ws.Range("N132").Value = -13885.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 75 (hunk 34)
$ws.Range("H75").Value = 40600
$ws.Range("J75").Value = 40600
$ws.Range("L75").Value = 40600
$ws.Range("N75").Value = -42472

# Row 78 (hunk 35)
$ws.Range("H78").Value = 40600
$ws.Range("J78").Value = 40600
$ws.Range("L78").Value = 121800
$ws.Range("N78").Value = -131160

# Row 92 (hunk 36)
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 104 (hunk 37)
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

# Row 105 (hunk 38)
$ws.Range("H105").Value = 43380
$ws.Range("J105").Value = 43380
$ws.Range("L105").Value = 43380
$ws.Range("N105").Value = -50368

# Row 132 (hunk 39)
$ws.Range("H132").Value = 4006.54
$ws.Range("I132").Value = 1951.8292
$ws.Range("J132").Value = 13366.889
$ws.Range("K132").Value = 5855.487599999999
$ws.Range("L132").Value = 40100.667
$ws.Range("M132").Value = -3325.487599999999
$ws.Range("N132").Value = -45160.667
